$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at 21, pushing the old rows 21+ (incl. the
#    signature block originally on rows 25/26) down by one, to rows 26/27.
$ws.Rows("21:21").Insert()

# 2) Row 20 still holds the "last row" heavy-bottom-border formatting
#    (periodo 2507). Copy that formatting onto the freshly inserted,
#    still-blank row 21, so row 21 becomes the new "last" row.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Turn row 20 back into a normal/middle table row by copying the
#    formatting from row 19 (a normal row) onto it.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Fill in row 21 with the new worker/period entry (same worker data,
#    new period 2508).
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143392800"
$ws.Range("D21").Value = "JORGE ENRIQUE CIRO TORO"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 132000
$ws.Range("G21").Value = 3300000

# 5) Update the two changed totals.
$ws.Range("E11").Value = 792000
$ws.Range("F13").Value = 6
